$wb = $excel.ActiveWorkbook

# Update the form_version value on the "settings" sheet.
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B3").Value = 20130408

# Make "settings" the active sheet/tab, with B3 selected (matches the
# saved selection state after editing the version value).
$settings.Activate()
$settings.Range("B3").Select()
